$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.761.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "'2.618.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").Value = "'515.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "'154.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "'2.632.18"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'6.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.16%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").Value = "'3.075.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'60.725.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "'21.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "'2.624.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'358.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.25%  "
$ws.Range("D21").Value = "'10.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").Value = "'6.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'60.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").Value = "'2.733.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "'7.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "'19.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "'5.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.74%  "
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'151.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'0.889"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.52%  "
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").Value = "'36.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'294.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.38%  "
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "'0.624"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "'0.0557"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "'19.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").Value = "'4.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").Value = "'10.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
